# Generate Report for Handoff
#
# A new handoff just completed for file
# 36c20b26-eb87-49fe-9993-e7bf92277144.md, so its "latest handoff"
# timestamps need to be refreshed on the Overview sheet and on each
# per-locale detail sheet.

$wb = $excel.ActiveWorkbook

$fileId = "36c20b26-eb87-49fe-9993-e7bf92277144"

# --- Overview sheet: column D = "Latest Handoff Date" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$overviewRows = $wsOverview.UsedRange.Rows.Count
for ($r = 1; $r -le $overviewRows; $r++) {
    $name = $wsOverview.Cells.Item($r, 1).Text
    if ($name -like "$fileId*") {
        $wsOverview.Cells.Item($r, 4).Value = "2016-03-19 07:29:12"
        break
    }
}

# --- Locale detail sheets: column E = "Latest Handoff Datetime" ---
$localeStamps = @{ "zh-cn" = "2016-03-19 07:29:04"; "de-de" = "2016-03-19 07:29:12" }

foreach ($localeName in $localeStamps.Keys) {
    $ws = $wb.Worksheets.Item($localeName)
    $rows = $ws.UsedRange.Rows.Count
    for ($r = 1; $r -le $rows; $r++) {
        $name = $ws.Cells.Item($r, 1).Text
        if ($name -like "$fileId*") {
            $ws.Cells.Item($r, 5).Value = $localeStamps[$localeName]
            break
        }
    }
}
